# Updating plots to incorporate Courtney's suggestions:
#  - Rename the "n" header (B1) and "count" header (G1) to "# species"
#  - Populate the new "# species" column G (rows 2-7) with the same
#    species counts already present in column B
#  - Leave the final selection on H14, matching the last-saved cursor spot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers
$ws.Range("B1").Value = "# species"
$ws.Range("G1").Value = "# species"

# Fill in the new "# species" values in column G, mirroring column B
# (use Value2 to read back plain numbers from column B)
$ws.Range("G2").Value = $ws.Range("B2").Value2
$ws.Range("G3").Value = $ws.Range("B3").Value2
$ws.Range("G4").Value = $ws.Range("B4").Value2
$ws.Range("G5").Value = $ws.Range("B5").Value2
$ws.Range("G6").Value = $ws.Range("B6").Value2
$ws.Range("G7").Value = $ws.Range("B7").Value2

# Match final cell selection recorded in the saved workbook
$ws.Range("H14").Select()
